$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 65.22210184182016
$ws.Range("C2").Value = 64.18251735683648
$ws.Range("D2").Value = 65.22210184182016
$ws.Range("E2").Value = 64.51749665215661
$ws.Range("F2").Value = 0.7028552813122796
$ws.Range("G2").Value = 35.05518933006285
